$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.590.57'
$ws.Range('E2').Value = '  +2.43%  '
$ws.Range('D3').Value = '3.362.92'
$ws.Range('E3').Value = '  +2.88%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = '''191.41'
$ws.Range('E5').Value = '  +2.84%  '
$ws.Range('D6').Value = '''591.52'
$ws.Range('E6').Value = '  +1.76%  '
$ws.Range('B7').Value = 'USDC'
$ws.Range('C7').Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range('D7').Value = '''1.00'
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('B8').Value = 'XRP'
$ws.Range('C8').Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range('D8').Value = '''0.609'
$ws.Range('E8').Value = '  +1.18%  '
$ws.Range('E9').Value = '  +1.83%  '
$ws.Range('D10').Value = '''6.77'
$ws.Range('E10').Value = '  +2.92%  '
$ws.Range('D11').Value = '''0.417'
$ws.Range('E11').Value = '  +1.43%  '
$ws.Range('D12').Value = '3.953.13'
$ws.Range('E12').Value = '  +3.03%  '
$ws.Range('E13').Value = '  -1.08%  '
$ws.Range('D14').Value = '''28.51'
$ws.Range('E14').Value = '  +3.44%  '
$ws.Range('D15').Value = '69.619.21'
$ws.Range('E15').Value = '  +2.45%  '
$ws.Range('E16').Value = '  +1.19%  '
$ws.Range('D17').Value = '3.375.65'
$ws.Range('E17').Value = '  +3.46%  '
$ws.Range('D18').Value = '''453.52'
$ws.Range('E18').Value = '  +13.91%  '
$ws.Range('D19').Value = '''5.81'
$ws.Range('E19').Value = '  +0.98%  '
$ws.Range('D20').Value = '''13.72'
$ws.Range('E20').Value = '  +1.23%  '
$ws.Range('D21').Value = '''7.86'
$ws.Range('E21').Value = '  +2.99%  '
$ws.Range('D22').Value = '''75.92'
$ws.Range('E22').Value = '  +6.12%  '
$ws.Range('D23').Value = '''0.999'
$ws.Range('E23').Value = '  -0.10%  '
$ws.Range('D24').Value = '3.522.34'
$ws.Range('E24').Value = '  +3.28%  '
$ws.Range('D25').Value = '''0.522'
$ws.Range('E25').Value = '  +1.90%  '
$ws.Range('E26').Value = '  +3.06%  '
$ws.Range('E27').Value = '  +1.69%  '
$ws.Range('D28').Value = '''9.44'
$ws.Range('E28').Value = '  -0.86%  '
$ws.Range('D29').Value = '''1.00'
$ws.Range('E29').Value = '  -0.54%  '
$ws.Range('D30').Value = '''2.01'
$ws.Range('E30').Value = '  +2.78%  '
$ws.Range('D31').Value = '''23.30'
$ws.Range('E31').Value = '  +2.61%  '
$ws.Range('D32').Value = '''5.54'
$ws.Range('E32').Value = '  +0.63%  '
$ws.Range('E33').Value = '  +2.35%  '
$ws.Range('D34').Value = '''6.96'
$ws.Range('E34').Value = '  -0.06%  '
$ws.Range('D35').Value = '''0.998'
$ws.Range('E35').Value = '  -0.02%  '
$ws.Range('E36').Value = '  +5.95%  '
$ws.Range('D37').Value = '''164.57'
$ws.Range('E37').Value = '  +0.64%  '
$ws.Range('D38').Value = '''1.94'
$ws.Range('E38').Value = '  +2.28%  '
$ws.Range('D39').Value = '''27.16'
$ws.Range('E39').Value = '  +1.51%  '
$ws.Range('D40').Value = '''0.810'
$ws.Range('E40').Value = '  +0.00%  '
$ws.Range('D41').Value = '''4.59'
$ws.Range('E41').Value = '  +1.09%  '
$ws.Range('E42').Value = '  +1.60%  '
$ws.Range('D43').Value = '2.721.41'
$ws.Range('E43').Value = '  +1.86%  '
$ws.Range('D44').Value = '''2.51'
$ws.Range('E44').Value = '  +2.98%  '
$ws.Range('D45').Value = '''0.0688'
$ws.Range('E45').Value = '  +0.15%  '
$ws.Range('D46').Value = '''25.36'
$ws.Range('E46').Value = '  +2.12%  '
$ws.Range('D47').Value = '''40.92'
$ws.Range('E47').Value = '  +0.44%  '
$ws.Range('D48').Value = '''334.99'
$ws.Range('E48').Value = '  +0.45%  '
$ws.Range('D49').Value = '''0.0284'
$ws.Range('E49').Value = '  +2.44%  '
$ws.Range('D50').Value = '''32.50'
$ws.Range('E50').Value = '  +5.33%  '
$ws.Range('E51').Value = '  +3.90%  '
